$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated monthly_mean (column D) and weekly_share (column F) values,
# recalculated after dividing monthly_mean by the number of weeks in
# each month group (weekly_mean in column E stays the same).

$updates = @(
    @{ Row = 2;  D = 70.76300000000001;  F = 1.127 },
    @{ Row = 3;  D = 70.76300000000001;  F = 1.004 },
    @{ Row = 4;  D = 70.76300000000001;  F = 1.347 },
    @{ Row = 5;  D = 70.76300000000001;  F = 0.651 },
    @{ Row = 6;  D = 70.76300000000001;  F = 0.871 },
    @{ Row = 7;  D = 74.67700000000001;  F = 1.159 },
    @{ Row = 8;  D = 74.67700000000001;  F = 0.902 },
    @{ Row = 9;  D = 74.67700000000001;  F = 0.788 },
    @{ Row = 10; D = 74.67700000000001;  F = 1.151 },
    @{ Row = 11; D = 70.078;             F = 1.346 },
    @{ Row = 12; D = 70.078;             F = 1.133 },
    @{ Row = 13; D = 70.078;             F = 0.842 },
    @{ Row = 14; D = 70.078;             F = 0.679 },
    @{ Row = 16; D = 34.706;             F = 1.095 },
    @{ Row = 17; D = 34.706;             F = 0.905 },
    @{ Row = 18; D = 23.003;             F = 1.112 },
    @{ Row = 19; D = 23.003;             F = 1.052 },
    @{ Row = 20; D = 23.003;             F = 0.991 },
    @{ Row = 21; D = 23.003;             F = 0.845 },
    @{ Row = 22; D = 16.749;             F = 1.042 },
    @{ Row = 23; D = 16.749;             F = 0.959 },
    @{ Row = 24; D = 16.749;             F = 1.046 },
    @{ Row = 25; D = 16.749;             F = 0.953 },
    @{ Row = 26; D = 10.489;             F = 1.204 },
    @{ Row = 27; D = 10.489;             F = 0.796 },
    @{ Row = 28; D = 14.778;             F = 0.866 },
    @{ Row = 29; D = 14.778;             F = 1.134 },
    @{ Row = 31; D = 89.745;             F = 0.872 },
    @{ Row = 32; D = 89.745;             F = 1.268 },
    @{ Row = 33; D = 89.745;             F = 1.027 },
    @{ Row = 34; D = 89.745;             F = 0.833 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}
